$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("E4").Value = 0.485
$ws.Range("N4").Value = 0.48
$ws.Range("O4").Value = 0.065
$ws.Range("P4").Value = 0.255
$ws.Range("Q4").Value = 0.055
$ws.Range("R4").Value = 0.037
$ws.Range("S4").Value = 0.192
$ws.Range("W4").Value = 0.385
$ws.Range("X4").Value = 0.106
$ws.Range("Y4").Value = 0.326
$ws.Range("AI4").Value = 0.402
$ws.Range("AJ4").Value = 0.094
$ws.Range("AK4").Value = 0.307
$ws.Range("AU4").Value = 0.257
$ws.Range("AV4").Value = 0.023
$ws.Range("AW4").Value = 0.152
$ws.Range("BA4").Value = 2.006
$ws.Range("BB4").Value = 0.145
$ws.Range("BC4").Value = 0.381
$ws.Range("BG4").Value = 0.708
$ws.Range("BH4").Value = 0.145
$ws.Range("BI4").Value = 0.381
$ws.Range("BM4").Value = 0.738
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.252
$ws.Range("BP4").Value = 0.669
$ws.Range("BQ4").Value = 0.755
$ws.Range("E5").Value = 0.602
$ws.Range("F5").Value = 0.062
$ws.Range("G5").Value = 0.25
$ws.Range("N5").Value = 0.715
$ws.Range("O5").Value = 0.078
$ws.Range("P5").Value = 0.279
$ws.Range("Q5").Value = 0.036
$ws.Range("R5").Value = 0.015
$ws.Range("S5").Value = 0.123
$ws.Range("W5").Value = 0.348
$ws.Range("X5").Value = 0.098
$ws.Range("Y5").Value = 0.313
$ws.Range("AI5").Value = 0.395
$ws.Range("AJ5").Value = 0.083
$ws.Range("AK5").Value = 0.288
$ws.Range("AU5").Value = 0.478
$ws.Range("AV5").Value = 0.067
$ws.Range("AW5").Value = 0.259
$ws.Range("BA5").Value = 1.289
$ws.Range("BB5").Value = 0.077
$ws.Range("BC5").Value = 0.277
$ws.Range("BG5").Value = 0.369
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.225
$ws.Range("BM5").Value = 0.527
$ws.Range("BN5").Value = 0.049
$ws.Range("BO5").Value = 0.222
$ws.Range("BP5").Value = 0.43
$ws.Range("BQ5").Value = 0.455
$ws.Range("E6").Value = 0.537
$ws.Range("N6").Value = 0.574
$ws.Range("Q6").Value = 0.044
$ws.Range("W6").Value = 0.366
$ws.Range("AI6").Value = 0.398
$ws.Range("AU6").Value = 0.334
$ws.Range("BA6").Value = 1.562
$ws.Range("BG6").Value = 0.485
$ws.Range("BM6").Value = 0.615
$ws.Range("BP6").Value = 0.521
$ws.Range("BQ6").Value = 0.564
$ws.Range("E7").Value = 0.574
$ws.Range("N7").Value = 0.651
$ws.Range("Q7").Value = 0.039
$ws.Range("W7").Value = 0.355
$ws.Range("AI7").Value = 0.396
$ws.Range("AU7").Value = 0.408
$ws.Range("BA7").Value = 1.385
$ws.Range("BG7").Value = 0.408
$ws.Range("BM7").Value = 0.559
$ws.Range("BP7").Value = 0.462
$ws.Range("BQ7").Value = 0.492
$ws.Range("E8").Value = 0.698
$ws.Range("F8").Value = 0.075
$ws.Range("G8").Value = 0.273
$ws.Range("N8").Value = 0.797
$ws.Range("O8").Value = 0.062
$ws.Range("P8").Value = 0.248
$ws.Range("Q8").Value = 0.04
$ws.Range("S8").Value = 0.161
$ws.Range("W8").Value = 0.414
$ws.Range("X8").Value = 0.115
$ws.Range("Y8").Value = 0.339
$ws.Range("AI8").Value = 0.462
$ws.Range("AJ8").Value = 0.132
$ws.Range("AK8").Value = 0.363
$ws.Range("AU8").Value = 0.423
$ws.Range("AV8").Value = 0.077
$ws.Range("AW8").Value = 0.277
$ws.Range("BA8").Value = 1.733
$ws.Range("BB8").Value = 0.11
$ws.Range("BC8").Value = 0.331
$ws.Range("BG8").Value = 0.545
$ws.Range("BH8").Value = 0.111
$ws.Range("BI8").Value = 0.334
$ws.Range("BM8").Value = 0.676
$ws.Range("BN8").Value = 0.061
$ws.Range("BO8").Value = 0.247
$ws.Range("BP8").Value = 0.578
$ws.Range("BQ8").Value = 0.618
$ws.Range("E9").Value = 0.65
$ws.Range("F9").Value = 0.227
$ws.Range("G9").Value = 0.477
$ws.Range("N9").Value = 0.725
$ws.Range("O9").Value = 0.199
$ws.Range("P9").Value = 0.447
$ws.Range("W9").Value = 0.3
$ws.Range("X9").Value = 0.21
$ws.Range("Y9").Value = 0.458
$ws.Range("AI9").Value = 0.4
$ws.Range("AJ9").Value = 0.24
$ws.Range("AK9").Value = 0.49
$ws.Range("BA9").Value = 1.675
$ws.Range("BG9").Value = 0.575
$ws.Range("BH9").Value = 0.244
$ws.Range("BI9").Value = 0.494
$ws.Range("BM9").Value = 0.65
$ws.Range("BN9").Value = 0.227
$ws.Range("BO9").Value = 0.477
$ws.Range("BP9").Value = 0.558
$ws.Range("BQ9").Value = 0.604
$ws.Range("E10").Value = 0.8
$ws.Range("F10").Value = 0.16
$ws.Range("G10").Value = 0.4
$ws.Range("N10").Value = 0.925
$ws.Range("O10").Value = 0.069
$ws.Range("P10").Value = 0.263
$ws.Range("W10").Value = 0.525
$ws.Range("X10").Value = 0.249
$ws.Range("Y10").Value = 0.499
$ws.Range("AI10").Value = 0.5
$ws.Range("AJ10").Value = 0.25
$ws.Range("AK10").Value = 0.5
$ws.Range("AU10").Value = 0.425
$ws.Range("AV10").Value = 0.244
$ws.Range("AW10").Value = 0.494
$ws.Range("BA10").Value = 2.15
$ws.Range("BB10").Value = 0.219
$ws.Range("BC10").Value = 0.468
$ws.Range("BG10").Value = 0.65
$ws.Range("BH10").Value = 0.227
$ws.Range("BI10").Value = 0.477
$ws.Range("BM10").Value = 0.825
$ws.Range("BN10").Value = 0.144
$ws.Range("BO10").Value = 0.38
$ws.Range("BP10").Value = 0.717
$ws.Range("BQ10").Value = 0.755
$ws.Range("E11").Value = 0.85
$ws.Range("F11").Value = 0.127
$ws.Range("G11").Value = 0.357
$ws.Range("N11").Value = 0.925
$ws.Range("O11").Value = 0.069
$ws.Range("P11").Value = 0.263
$ws.Range("W11").Value = 0.525
$ws.Range("X11").Value = 0.249
$ws.Range("Y11").Value = 0.499
$ws.Range("AI11").Value = 0.575
$ws.Range("AJ11").Value = 0.244
$ws.Range("AK11").Value = 0.494
$ws.Range("AU11").Value = 0.6
$ws.Range("AV11").Value = 0.24
$ws.Range("AW11").Value = 0.49
$ws.Range("BA11").Value = 2.15
$ws.Range("BB11").Value = 0.219
$ws.Range("BC11").Value = 0.468
$ws.Range("BG11").Value = 0.65
$ws.Range("BH11").Value = 0.227
$ws.Range("BI11").Value = 0.477
$ws.Range("BM11").Value = 0.825
$ws.Range("BN11").Value = 0.144
$ws.Range("BO11").Value = 0.38
$ws.Range("BP11").Value = 0.717
$ws.Range("BQ11").Value = 0.761
$ws.Range("E12").Value = 1.441
$ws.Range("F12").Value = 0.894
$ws.Range("G12").Value = 0.945
$ws.Range("N12").Value = 1.27
$ws.Range("O12").Value = 0.305
$ws.Range("P12").Value = 0.553
$ws.Range("W12").Value = 1.524
$ws.Range("X12").Value = 0.44
$ws.Range("Y12").Value = 0.663
$ws.Range("AI12").Value = 1.652
$ws.Range("AJ12").Value = 1.531
$ws.Range("AK12").Value = 1.237
$ws.Range("AU12").Value = 2.846
$ws.Range("AV12").Value = 3.361
$ws.Range("AW12").Value = 1.833
$ws.Range("BA12").Value = 3.84
$ws.Range("BB12").Value = 0.469
$ws.Range("BC12").Value = 0.685
$ws.Range("BG12").Value = 1.154
$ws.Range("BH12").Value = 0.207
$ws.Range("BI12").Value = 0.455
$ws.Range("BM12").Value = 1.242
$ws.Range("BN12").Value = 0.244
$ws.Range("BO12").Value = 0.494
$ws.Range("BP12").Value = 1.28
$ws.Range("BQ12").Value = 1.26
$ws.Range("E13").Value = 1.373
$ws.Range("F13").Value = 0.27
$ws.Range("G13").Value = 0.52
$ws.Range("N13").Value = 1.729
$ws.Range("O13").Value = 0.497
$ws.Range("P13").Value = 0.705
$ws.Range("W13").Value = 0.958
$ws.Range("X13").Value = 0.19
$ws.Range("Y13").Value = 0.435
$ws.Range("AI13").Value = 1.148
$ws.Range("AJ13").Value = 0.317
$ws.Range("AK13").Value = 0.563
$ws.Range("AU13").Value = 2.039
$ws.Range("AV13").Value = 0.339
$ws.Range("AW13").Value = 0.582
$ws.Range("BA13").Value = 2.204
$ws.Range("BB13").Value = 0.291
$ws.Range("BC13").Value = 0.54
$ws.Range("BG13").Value = 0.541
$ws.Range("BH13").Value = 0.053
$ws.Range("BI13").Value = 0.231
$ws.Range("BM13").Value = 0.801
$ws.Range("BN13").Value = 0.166
$ws.Range("BO13").Value = 0.408
$ws.Range("BP13").Value = 0.735
$ws.Range("BQ13").Value = 0.669
